# The sheet's "summary" block is restructured:
#  - the per-branch line items (New/Carryover/Confirmed/Unconfirmed/Withdrawn)
#    get their branch name prefixed onto the label (e.g. "     New nominations"
#    becomes "     Civilian, New nominations") so each label is unique instead
#    of being reused across branches,
#  - the trailing totals block drops the old "Summary" / carried-over / received
#    labels in favor of "Total new nominations" / "Total carryover nominations"
#    and shifts the other totals up by one row,
#  - the final row (old row 44, "Total returned to the White House") is merged
#    up into what is now row 43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old last row; row 43's old content goes away and every row below it
# (none) would shift up - this leaves exactly 43 data rows (A1:B43), matching
# the new dimension.
$ws.Rows("44").Delete() | Out-Null

# --- Column A: prefix each branch's line-item labels with the branch name ---
$ws.Range("A7").Value = "     Civilian, New nominations"
$ws.Range("A8").Value = "     Civilian, Carryover nominations"
$ws.Range("A9").Value = "     Civilian, Confirmed "
$ws.Range("A10").Value = "     Civilian, Unconfirmed "
$ws.Range("A11").Value = "     Civilian, Withdrawn "

$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Carryover nominations"
$ws.Range("A15").Value = "     Other Civilian, Confirmed "
$ws.Range("A16").Value = "     Other Civilian, Unconfirmed "
$ws.Range("A17").Value = "     Other Civilian, Withdrawn "

$ws.Range("A19").Value = "     Air Force, New nominations"
$ws.Range("A20").Value = "     Air Force, Carryover nominations"
$ws.Range("A21").Value = "     Air Force, Confirmed "
$ws.Range("A22").Value = "     Air Force, Unconfirmed "
$ws.Range("A23").Value = "     Air Force, Withdrawn "

$ws.Range("A25").Value = "     Army, New nominations"
$ws.Range("A26").Value = "     Army, Carryover nominations"
$ws.Range("A27").Value = "     Army, Confirmed "
$ws.Range("A28").Value = "     Army, Unconfirmed "

$ws.Range("A30").Value = "     Navy, New nominations"
$ws.Range("A31").Value = "     Navy, Carryover nominations"
$ws.Range("A32").Value = "     Navy, Confirmed "
$ws.Range("A33").Value = "     Navy, Unconfirmed "

$ws.Range("A35").Value = "     Marine Corps, New nominations"
$ws.Range("A36").Value = "     Marine Corps, Confirmed "
$ws.Range("A37").Value = "     Marine Corps, Unconfirmed "

# --- Column A: rework the totals block labels (rows 38-43) ---
$ws.Range("A38").Value = "Total new nominations"
$ws.Range("A39").Value = "Total carryover nominations"
$ws.Range("A40").Value = "Total confirmed "
$ws.Range("A41").Value = "Total unconfirmed "
$ws.Range("A42").Value = "Total withdrawn "
$ws.Range("A43").Value = "Total returned to the White House "

# --- Column B: totals block values now start at row 38 ---
# B38 was previously empty (its row only held the "Summary" header) - give it
# the "#,##0"-formatted numeric style used by its sibling total cells.
$ws.Range("B38").Value = 21855
$ws.Range("B38").NumberFormat = "#,##0"

$ws.Range("B40").Value = 21751

# B41 changes from the "#,##0" style to the plain right-aligned numeric style;
# copy that formatting over from a cell that already has it, then set the value.
$ws.Range("B42").Copy() | Out-Null
$ws.Range("B41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B41").Value = 225

$ws.Range("B42").Value = 24
$ws.Range("B43").Value = 0
